$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "11÷4="
$t.Cell(1, 2).Range.Text = "34÷4="
$t.Cell(1, 3).Range.Text = "76÷2="
$t.Cell(1, 4).Range.Text = "64÷7="
$t.Cell(1, 5).Range.Text = "10÷5="

$t.Cell(5, 1).Range.Text = "76÷6="
$t.Cell(5, 2).Range.Text = "48÷8="
$t.Cell(5, 3).Range.Text = "98÷4="
$t.Cell(5, 4).Range.Text = "89÷4="
$t.Cell(5, 5).Range.Text = "83÷4="

$t.Cell(9, 1).Range.Text = "23÷9="
$t.Cell(9, 2).Range.Text = "78÷2="
$t.Cell(9, 3).Range.Text = "83÷8="
$t.Cell(9, 4).Range.Text = "52÷6="
$t.Cell(9, 5).Range.Text = "83÷8="

$t.Cell(13, 1).Range.Text = "81÷6="
$t.Cell(13, 2).Range.Text = "91÷6="
$t.Cell(13, 3).Range.Text = "43÷7="
$t.Cell(13, 4).Range.Text = "45÷3="
$t.Cell(13, 5).Range.Text = "19÷8="

$t.Cell(17, 1).Range.Text = "50÷2="
$t.Cell(17, 2).Range.Text = "25÷5="
$t.Cell(17, 3).Range.Text = "44÷9="
$t.Cell(17, 4).Range.Text = "29÷9="
$t.Cell(17, 5).Range.Text = "70÷7="
